# changement revue de code
#
# "Fiche de revue de code" tracking sheet: the reviewer re-triaged items
# #9-#11. The "bloc conditionnel if __name__" remark moves up to item #9,
# the "ligne vide a la fin du code" remark moves up to item #10, and item
# #11 becomes a fresh "fermer la base de donnees" remark whose LIGNE (col B)
# hasn't been filled in yet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 20 / item 9 ---------------------------------------------------
$ws.Range("B20").Value = 232
$ws.Range("D20").Value = 'Il manque un bloc conditionnel  if __name__ == "__main__": programme_principal()'

# --- Row 21 / item 10 ---------------------------------------------------
$ws.Range("C21").Value = "cosmétique"
$ws.Range("D21").Value = "Il est toujours préférable de laisse une ligne vide à la fin du code"

# --- Row 22 / item 11 ---------------------------------------------------
$ws.Range("B22").ClearContents()
$ws.Range("C22").Value = "erreur"
$ws.Range("D22").Value = "Il est important de fermer la base de données après avoir terminé de l'utiliser"

# --- Leave the sheet with the cursor where the reviewer left it --------
$ws.Activate()
$ws.Range("D26").Select() | Out-Null
